$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 513 (existing rows 513:538 shift down to 515:540)
$ws.Rows.Item(513).Insert()
$ws.Rows.Item(513).Insert()

# Row 513
$ws.Range("A513").Value = 5
$ws.Range("B513").Value = "Macroferia Regional de Talca"
$ws.Range("C513").Value = "Maule"
$ws.Range("D513").Value = 45041
$ws.Range("E513").Value = 7
$ws.Range("F513").Value = 100114014
$ws.Range("G513").Value = "Betarraga"
$ws.Range("H513").Value = "Sin especificar"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 3000
$ws.Range("K513").Value = 700
$ws.Range("L513").Value = 700
$ws.Range("M513").Value = 700
$ws.Range("N513").Value = "`$/paquete 5 unidades"
$ws.Range("O513").Value = "Región del Maule"
$ws.Range("P513").Value = 140
$ws.Range("Q513").Value = 5
$ws.Range("R513").Value = "Hortaliza"

# Row 514
$ws.Range("A514").Value = 5
$ws.Range("B514").Value = "Macroferia Regional de Talca"
$ws.Range("C514").Value = "Maule"
$ws.Range("D514").Value = 45041
$ws.Range("E514").Value = 7
$ws.Range("F514").Value = 100114014
$ws.Range("G514").Value = "Betarraga"
$ws.Range("H514").Value = "Sin especificar"
$ws.Range("I514").Value = "Segunda"
$ws.Range("J514").Value = 2000
$ws.Range("K514").Value = 600
$ws.Range("L514").Value = 600
$ws.Range("M514").Value = 600
$ws.Range("N514").Value = "`$/paquete 5 unidades"
$ws.Range("O514").Value = "Región del Maule"
$ws.Range("P514").Value = 120
$ws.Range("Q514").Value = 5
$ws.Range("R514").Value = "Hortaliza"
